$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*The API documentation is in documentation>Backend API documentation*") {
        $target = $p
    }
}

$idx = $target.Index
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($idx + 1)
$newPara.Range.Text = "Optionally run the “package.py” script in the util folder to install all python dependencies for the OCR script automatically."
